$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 3
    3 = 3
    4 = 1
    5 = 1
    6 = 1
    7 = 1
    8 = 2
    9 = 1
    10 = 1
    11 = 0
    12 = 0
    13 = 0
    14 = 1
    15 = 0
    16 = 1
    17 = 2
    18 = 0
    19 = 0
    20 = 0
    22 = 0
    23 = 0
    24 = 3
    25 = 0
    26 = 0
    27 = 1
    28 = 0
    29 = 1
    30 = 1
    31 = 0
    32 = 0
    33 = 1
    34 = 1
    35 = 2
    36 = 1
    37 = 1
    38 = 2
    39 = 2
    40 = 1
    41 = 1
    42 = 2
    43 = 2
    44 = 1
    45 = 1
    46 = 3
    47 = 1
    48 = 2
    49 = 1
    50 = 2
    51 = 4
    52 = 0
    53 = 3
    54 = 2
    55 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
